# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Update the computed K values for each game row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G)
$kValues = @{
    2  = 0
    3  = 2
    4  = 3
    5  = 3
    6  = 2
    7  = 1
    8  = 1
    9  = 0
    10 = 2
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    17 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
